{"js": "// Remove the stray \"<lb/>\" marker (rendered as two separate runs: \"<\" and\n// \"lb/>\") that immediately follows \"dans l'<m>eau</m>\" and immediately\n// precedes \"</ab>\" in the document body. Everything else - including the\n// formatting of the surrounding runs - must stay untouched.\nconst body = context.document.body;\n\n// Scope the search to the single paragraph that contains this text so we\n// don't touch any of the many other \"<lb/>\" line-break markers elsewhere\n// in the document.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.indexOf(\"dans l'\") !== -1\n);\n\nif (target) {\n  const matches = target.search(\"<lb/>\", { matchCase: true });\n  matches.load(\"items,text\");\n  await context.sync();\n\n  if (matches.items.length > 0) {\n    // Delete just the matched text range; this removes the underlying\n    // \"<\" run and \"lb/>\" run without merging/reformatting neighboring runs.\n    matches.items[0].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the stray \"<lb/>\" marker (stored as two separate runs: \"<\" and\n# \"lb/>\") that sits right after \"dans l'<m>eau</m>\" and right before\n# \"</ab>\" in the document body. Nothing else in the document - including\n# the formatting of the neighboring runs - should change.\n\n$d = $word.ActiveDocument\n\n# Locate the one paragraph containing this text so the edit can't touch\n# any of the many other \"<lb/>\" line-break markers elsewhere in the doc.\n$target = $null\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text -like \"*dans l'*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $r = $target.Range\n    $found = $r.Find.Execute(\"<lb/>\", $true)\n    if ($found) {\n        # Deleting the matched range removes exactly the \"<\" run and the\n        # \"lb/>\" run, leaving every other run (and its formatting) intact.\n        $r.Delete()\n    }\n}\n"}
